# Add a set of new HTML-course exercises into the tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Column K (exercise index within the "links" block) and column M/N (exercise
# names) for rows 20-27 — a new sub-section of tasks about the <a> tag.
$ws.Range("K20").Value = 1
$ws.Range("M20").Value = "Доктайп + html"
$ws.Range("M20").Style = $ws.Range("M9").Style

$ws.Range("K21").Value = 2
$ws.Range("M21").Value = "head + title"
$ws.Range("M21").Style = $ws.Range("M9").Style

$ws.Range("K22").Value = 3
$ws.Range("M22").Value = "body"
$ws.Range("M22").Style = $ws.Range("M9").Style

$ws.Range("K23").Value = 4
$ws.Range("M23").Value = "links"
$ws.Range("N23").Value = "external"

$ws.Range("K24").Value = 5
$ws.Range("M24").Value = "links"

$ws.Range("K25").Value = 6
$ws.Range("M25").Value = "links"
$ws.Range("N25").Value = "relative"

$ws.Range("N24").Value = "обернуть img"

$ws.Range("K26").Value = 7
$ws.Range("M26").Value = "links"

$ws.Range("K27").Value = 8
$ws.Range("M27").Value = "комментарии"

$ws.Range("N26").Value = "якорь на странице + заглушка"

# Apply the same fill style used by the earlier block (M7:M17) to the
# previously-unstyled continuation cells so the whole column reads uniformly.
$ws.Range("M10:M17").Style = $ws.Range("M9").Style

$ws.Range("M25").Select()
